# live_trading_results.xlsx update
# Trade #84 closed at 2026-02-17 21:17:50 - unknown UNKNOWN +0.000%
#
# This script:
#  1) Updates the Summary sheet aggregate stats
#  2) Updates the Strategy Status row for MarketMaking
#  3) Closes an existing open trade (Trade #112 -> early_exit) on both the
#     "All Trades" sheet and the "MarketMaking" sheet
#  4) Appends a brand-new open trade (Trade #145) to both the "All Trades"
#     sheet and the "MarketMaking" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel silently
# re-interpreting strings that look like dates/times (e.g. "2026-02-17"
# or "21:17:43") as date/time serial numbers. We temporarily force the
# cell to a Text number format, assign the raw value, then restore the
# cell back to the default "Normal" style so no visible formatting
# change is left behind.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# =======================================================================
# 1) Summary sheet
# =======================================================================
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.08   # Current Capital
$summary.Range("B4").Value = 0.87      # Total P&L $
$summary.Range("B5").Value = 0.16      # Total P&L %
$summary.Range("B6").Value = 112       # Total Trades
$summary.Range("B8").Value = 43        # Losing Trades
$summary.Range("B9").Value = 44.64     # Win Rate %

# =======================================================================
# 2) Strategy Status sheet (MarketMaking row)
# =======================================================================
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.08     # Capital
$status.Range("D5").Value = 79         # Trades
$status.Range("E5").Value = 0.76       # P&L $
$status.Range("F5").Value = 1.08       # P&L %
$status.Range("G5").Value = 45.57      # Win Rate %

# =======================================================================
# 3) All Trades sheet
# =======================================================================
$allTrades = $wb.Worksheets.Item("All Trades")

# --- 3a) Close existing Trade #112 (row 113) ---
$allTrades.Cells.Item(113, 7).Value  = 0.11                # Exit Price
Set-TextValue $allTrades.Cells.Item(113, 8) "CLOSED"       # Status
$allTrades.Cells.Item(113, 9).Value  = -35.2941            # P&L %
$allTrades.Cells.Item(113, 10).Value = -0.06               # P&L $
$allTrades.Cells.Item(113, 11).Value = 101.08              # Capital After
Set-TextValue $allTrades.Cells.Item(113, 12) "early_exit"  # Exit Reason
$allTrades.Cells.Item(113, 13).Value = 0.14                # Duration (min)

# --- 3b) Append brand-new open Trade #145 (row 146) ---
$allTrades.Cells.Item(146, 1).Value = 145                          # Trade #
Set-TextValue $allTrades.Cells.Item(146, 2) "2026-02-17"            # Date
Set-TextValue $allTrades.Cells.Item(146, 3) "21:17:43"              # Time
Set-TextValue $allTrades.Cells.Item(146, 4) "MarketMaking"          # Strategy
Set-TextValue $allTrades.Cells.Item(146, 5) "UP"                    # Side
$allTrades.Cells.Item(146, 6).Value  = 0.17                         # Entry Price
$allTrades.Cells.Item(146, 7).Value  = ""                           # Exit Price
Set-TextValue $allTrades.Cells.Item(146, 8) "OPEN"                  # Status
$allTrades.Cells.Item(146, 9).Value  = 0                            # P&L %
$allTrades.Cells.Item(146, 10).Value = 0                            # P&L $
$allTrades.Cells.Item(146, 11).Value = 101.1396151053151            # Capital After
$allTrades.Cells.Item(146, 12).Value = ""                           # Exit Reason
$allTrades.Cells.Item(146, 13).Value = 0                            # Duration (min)
$allTrades.Cells.Item(146, 14).Value = 0                            # Entry Slippage (bps)
$allTrades.Cells.Item(146, 15).Value = 0                            # Exit Slippage (bps)
$allTrades.Cells.Item(146, 16).Value = 0.6                          # Confidence
Set-TextValue $allTrades.Cells.Item(146, 17) "Normal spread capture: 19600 bps"  # Entry Reason

# =======================================================================
# 4) MarketMaking sheet
# =======================================================================
$mm = $wb.Worksheets.Item("MarketMaking")

# --- 4a) Close the same trade locally (row 80) ---
$mm.Cells.Item(80, 7).Value  = 0.11                # Exit Price
Set-TextValue $mm.Cells.Item(80, 8) "CLOSED"       # Status
$mm.Cells.Item(80, 9).Value  = -35.2941            # P&L %
$mm.Cells.Item(80, 10).Value = -0.06               # P&L $
$mm.Cells.Item(80, 11).Value = 101.08              # Capital After
Set-TextValue $mm.Cells.Item(80, 16) "early_exit"  # Exit Reason
$mm.Cells.Item(80, 17).Value = 0.14                # Duration (min)

# --- 4b) Append the same brand-new open trade locally (row 113) ---
$mm.Cells.Item(113, 1).Value = 145                          # Trade #
Set-TextValue $mm.Cells.Item(113, 2) "2026-02-17"            # Date
Set-TextValue $mm.Cells.Item(113, 3) "21:17:43"              # Time
Set-TextValue $mm.Cells.Item(113, 4) "MarketMaking"          # Strategy
Set-TextValue $mm.Cells.Item(113, 5) "UP"                    # Side
$mm.Cells.Item(113, 6).Value  = 0.17                         # Entry Price
$mm.Cells.Item(113, 7).Value  = ""                           # Exit Price
Set-TextValue $mm.Cells.Item(113, 8) "OPEN"                  # Status
$mm.Cells.Item(113, 9).Value  = 0                            # P&L %
$mm.Cells.Item(113, 10).Value = 0                            # P&L $
$mm.Cells.Item(113, 11).Value = 101.1396151053151            # Capital After
$mm.Cells.Item(113, 12).Value = 0                            # Entry Slippage (bps)
$mm.Cells.Item(113, 13).Value = 0                            # Exit Slippage (bps)
$mm.Cells.Item(113, 14).Value = 0.6                          # Confidence
Set-TextValue $mm.Cells.Item(113, 15) "Normal spread capture: 19600 bps"  # Entry Reason
$mm.Cells.Item(113, 16).Value = ""                           # Exit Reason
$mm.Cells.Item(113, 17).Value = 0                            # Duration (min)
